# Update the cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Price column (D) holds plain-text numeric-looking strings (e.g. thousand
    # separators rendered with dots, fixed trailing zeros, etc). Pre-format the
    # cell as Text so Excel stores the literal string instead of silently
    # coercing it to a floating point number (which would lose formatting like
    # trailing zeros, e.g. "0.150" -> 0.15).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "98.911.18"
$ws.Range("E2").Value = "  +1.01%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.342.84"
$ws.Range("E3").Value = "  +6.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - Solana
Set-TextValue "D5" "259.25"
$ws.Range("E5").Value = "  +7.92%  "

# Row 6 - BNB
Set-TextValue "D6" "625.31"
$ws.Range("E6").Value = "  +2.48%  "

# Row 7 - XRP
Set-TextValue "D7" "1.46"
$ws.Range("E7").Value = "  +31.67%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  +2.52%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.08%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.907"
$ws.Range("E10").Value = "  +15.33%  "

# Row 11 - LidoStakedEther
Set-TextValue "D11" "3.342.07"
$ws.Range("E11").Value = "  +6.19%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.28%  "

# Row 13 - Avalanche
Set-TextValue "D13" "37.94"
$ws.Range("E13").Value = "  +11.37%  "

# Row 14 - WrappedBTC
Set-TextValue "D14" "98.646.59"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15 - ShibaInu
Set-TextValue "D15" "0.0000250"
$ws.Range("E15").Value = "  +4.08%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.961.39"
$ws.Range("E16").Value = "  +6.16%  "

# Row 17 - Toncoin
$ws.Range("E17").Value = "  +1.57%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.339.67"
$ws.Range("E18").Value = "  +6.17%  "

# Row 19 - SuiNetwork
Set-TextValue "D19" "3.57"
$ws.Range("E19").Value = "  +2.56%  "

# Row 20 - Chainlink
Set-TextValue "D20" "15.32"
$ws.Range("E20").Value = "  +4.76%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "492.35"
$ws.Range("E21").Value = "  -5.22%  "

# Row 22 - Polkadot
Set-TextValue "D22" "6.16"
$ws.Range("E22").Value = "  +7.66%  "

# Row 23 - PEPE
$ws.Range("E23").Value = "  +9.55%  "

# Row 24 - Uniswap
Set-TextValue "D24" "9.45"
$ws.Range("E24").Value = "  +7.08%  "

# Row 25 - NEARProtocol
Set-TextValue "D25" "5.66"
$ws.Range("E25").Value = "  +3.46%  "

# Row 26 - Litecoin
Set-TextValue "D26" "90.03"
$ws.Range("E26").Value = "  +1.70%  "

# Row 27 - Aptos
Set-TextValue "D27" "11.99"
$ws.Range("E27").Value = "  +3.12%  "

# Rows 28 & 29 swapped rank: Stellar now ranks above WrappedeETH
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D28" "0.298"
$ws.Range("E28").Value = "  +25.90%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D29" "3.517.34"
$ws.Range("E29").Value = "  +6.14%  "

# Row 30 - Dai
$ws.Range("E30").Value = "  -0.06%  "

# Row 31 - Cronos
Set-TextValue "D31" "0.191"
$ws.Range("E31").Value = "  +7.77%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.137"
$ws.Range("E32").Value = "  +13.26%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +8.36%  "

# Row 34 - Binance-PegBSC-USD
$ws.Range("E34").Value = "  +0.17%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "28.41"
$ws.Range("E35").Value = "  +6.38%  "

# Rows 36 & 37 swapped rank: Kaspa now ranks above RenderToken
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.150"
$ws.Range("E36").Value = "  -1.09%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D37" "7.29"
$ws.Range("E37").Value = "  +0.95%  "

# Row 38 - PancakeSwap
$ws.Range("E38").Value = "  +4.09%  "

# Row 39 - Bittensor
Set-TextValue "D39" "502.29"
$ws.Range("E39").Value = "  +7.76%  "

# Row 40 - PolygonEcosystemToken
Set-TextValue "D40" "0.462"
$ws.Range("E40").Value = "  +6.29%  "

# Row 41 - WhiteBITCoin
Set-TextValue "D41" "24.89"
$ws.Range("E41").Value = "  +2.14%  "

# Row 42 - Fetch.AI
$ws.Range("E42").Value = "  +3.91%  "

# Row 43 - MantraDAO
$ws.Range("E43").Value = "  +4.25%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  +7.51%  "

# Row 45 - USDe (unchanged)

# Row 46 - ARBITRUM
Set-TextValue "D46" "0.782"
$ws.Range("E46").Value = "  +11.57%  "

# Row 47 - Monero
Set-TextValue "D47" "159.33"
$ws.Range("E47").Value = "  -2.07%  "

# Row 48 - Stacks
Set-TextValue "D48" "1.96"
$ws.Range("E48").Value = "  +1.46%  "

# Row 49 - Mantle
$ws.Range("E49").Value = "  +7.66%  "

# Row 50 - Filecoin
Set-TextValue "D50" "4.70"
$ws.Range("E50").Value = "  +3.18%  "

# Row 51 - OKB
Set-TextValue "D51" "45.94"
$ws.Range("E51").Value = "  +4.42%  "

Write-Output "Applied cryptos update"
